# Auto-generated update for "Recommandations" and "Top_YTD" sheets
# Reflects the latest BRVM automated recommendation refresh (GitHub Actions).

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd  = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": refresh existing rows 2-41 and append new rows 42-44 ---
# Row 2: BRVM - SERVICES PUBLICS
$wsReco.Range("A2").Value = 'BRVM - SERVICES PUBLICS'
$wsReco.Range("B2").Value = 0
$wsReco.Range("C2").Value = 8
$wsReco.Range("D2").Value = 3221.59
$wsReco.Range("E2").Value = 100.18
$wsReco.Range("F2").Value = '🟡 Observer'
$wsReco.Range("G2").Value = '➖ Neutre'

# Row 3: BRVM - AUTRES SECTEURS
$wsReco.Range("A3").Value = 'BRVM - AUTRES SECTEURS'
$wsReco.Range("B3").Value = 0
$wsReco.Range("C3").Value = 4
$wsReco.Range("D3").Value = 2552.3
$wsReco.Range("E3").Value = 640.73
$wsReco.Range("F3").Value = '🟡 Observer'
$wsReco.Range("G3").Value = '➖ Neutre'

# Row 4: NEI-CEDA CI
$wsReco.Range("A4").Value = 'NEI-CEDA CI'
$wsReco.Range("B4").Value = 0
$wsReco.Range("C4").Value = 4
$wsReco.Range("D4").Value = 2385
$wsReco.Range("E4").Value = 595
$wsReco.Range("F4").Value = '🟡 Observer'
$wsReco.Range("G4").Value = '➖ Neutre'

# Row 5: UNIWAX CI
$wsReco.Range("A5").Value = 'UNIWAX CI'
$wsReco.Range("B5").Value = 0
$wsReco.Range("C5").Value = 4
$wsReco.Range("D5").Value = 2290
$wsReco.Range("E5").Value = 580
$wsReco.Range("F5").Value = '🟡 Observer'
$wsReco.Range("G5").Value = '➖ Neutre'

# Row 6: AIR LIQUIDE CI
$wsReco.Range("A6").Value = 'AIR LIQUIDE CI'
$wsReco.Range("B6").Value = 0
$wsReco.Range("C6").Value = 4
$wsReco.Range("D6").Value = 2180
$wsReco.Range("E6").Value = 550
$wsReco.Range("F6").Value = '🟡 Observer'
$wsReco.Range("G6").Value = '➖ Neutre'

# Row 7: SETAO CI
$wsReco.Range("A7").Value = 'SETAO CI'
$wsReco.Range("B7").Value = 0
$wsReco.Range("C7").Value = 4
$wsReco.Range("D7").Value = 2165
$wsReco.Range("E7").Value = 540
$wsReco.Range("F7").Value = '🟡 Observer'
$wsReco.Range("G7").Value = '➖ Neutre'

# Row 8: CFAO MOTORS CI
$wsReco.Range("A8").Value = 'CFAO MOTORS CI'
$wsReco.Range("B8").Value = 0
$wsReco.Range("C8").Value = 3
$wsReco.Range("D8").Value = 1915
$wsReco.Range("E8").Value = 635
$wsReco.Range("F8").Value = '🟡 Observer'
$wsReco.Range("G8").Value = '➖ Neutre'

# Row 9: BRVM - DISTRIBUTION
$wsReco.Range("A9").Value = 'BRVM - DISTRIBUTION'
$wsReco.Range("B9").Value = 0
$wsReco.Range("C9").Value = 4
$wsReco.Range("D9").Value = 1458.44
$wsReco.Range("E9").Value = 368.81
$wsReco.Range("F9").Value = '🟡 Observer'
$wsReco.Range("G9").Value = '➖ Neutre'

# Row 10: BRVM - TRANSPORT
$wsReco.Range("A10").Value = 'BRVM - TRANSPORT'
$wsReco.Range("B10").Value = 0
$wsReco.Range("C10").Value = 4
$wsReco.Range("D10").Value = 1395.21
$wsReco.Range("E10").Value = 348.8
$wsReco.Range("F10").Value = '🟡 Observer'
$wsReco.Range("G10").Value = '➖ Neutre'

# Row 11: BRVM - AGRICULTURE
$wsReco.Range("A11").Value = 'BRVM - AGRICULTURE'
$wsReco.Range("B11").Value = 0
$wsReco.Range("C11").Value = 4
$wsReco.Range("D11").Value = 1286.97
$wsReco.Range("E11").Value = 325.72
$wsReco.Range("F11").Value = '🟡 Observer'
$wsReco.Range("G11").Value = '➖ Neutre'

# Row 12: BRVM - INDUSTRIE
$wsReco.Range("A12").Value = 'BRVM - INDUSTRIE'
$wsReco.Range("B12").Value = 0
$wsReco.Range("C12").Value = 4
$wsReco.Range("D12").Value = 1054.49
$wsReco.Range("E12").Value = 263.45
$wsReco.Range("F12").Value = '🟡 Observer'
$wsReco.Range("G12").Value = '➖ Neutre'

# Row 13: BRVM - CONSOMMATION DE BASE
$wsReco.Range("A13").Value = 'BRVM - CONSOMMATION DE BASE'
$wsReco.Range("B13").Value = 0
$wsReco.Range("C13").Value = 4
$wsReco.Range("D13").Value = 869.95
$wsReco.Range("E13").Value = 217.89
$wsReco.Range("F13").Value = '🟡 Observer'
$wsReco.Range("G13").Value = '➖ Neutre'

# Row 14: BRVM-PRINCIPAL
$wsReco.Range("A14").Value = 'BRVM-PRINCIPAL'
$wsReco.Range("B14").Value = 0
$wsReco.Range("C14").Value = 4
$wsReco.Range("D14").Value = 756.89
$wsReco.Range("E14").Value = 190.5
$wsReco.Range("F14").Value = '🟡 Observer'
$wsReco.Range("G14").Value = '➖ Neutre'

# Row 15: BRVM - INDUSTRIELS
$wsReco.Range("A15").Value = 'BRVM - INDUSTRIELS'
$wsReco.Range("B15").Value = 0
$wsReco.Range("C15").Value = 4
$wsReco.Range("D15").Value = 547.18
$wsReco.Range("E15").Value = 136.99
$wsReco.Range("F15").Value = '🟡 Observer'
$wsReco.Range("G15").Value = '➖ Neutre'

# Row 16: BRVM-PRESTIGE
$wsReco.Range("A16").Value = 'BRVM-PRESTIGE'
$wsReco.Range("B16").Value = 0
$wsReco.Range("C16").Value = 4
$wsReco.Range("D16").Value = 514.52
$wsReco.Range("E16").Value = 128.82
$wsReco.Range("F16").Value = '🟡 Observer'
$wsReco.Range("G16").Value = '➖ Neutre'

# Row 17: BRVM - FINANCES
$wsReco.Range("A17").Value = 'BRVM - FINANCES'
$wsReco.Range("B17").Value = 0
$wsReco.Range("C17").Value = 4
$wsReco.Range("D17").Value = 485.98
$wsReco.Range("E17").Value = 122.29
$wsReco.Range("F17").Value = '🟡 Observer'
$wsReco.Range("G17").Value = '➖ Neutre'

# Row 18: BRVM - SERVICES FINANCIERS
$wsReco.Range("A18").Value = 'BRVM - SERVICES FINANCIERS'
$wsReco.Range("B18").Value = 0
$wsReco.Range("C18").Value = 4
$wsReco.Range("D18").Value = 477.61
$wsReco.Range("E18").Value = 120.19
$wsReco.Range("F18").Value = '🟡 Observer'
$wsReco.Range("G18").Value = '➖ Neutre'

# Row 19: BRVM - ENERGIE
$wsReco.Range("A19").Value = 'BRVM - ENERGIE'
$wsReco.Range("B19").Value = 0
$wsReco.Range("C19").Value = 4
$wsReco.Range("D19").Value = 435.51
$wsReco.Range("E19").Value = 110.28
$wsReco.Range("F19").Value = '🟡 Observer'
$wsReco.Range("G19").Value = '➖ Neutre'

# Row 20: BRVM - CONSOMMATION DISCRETIONNAIRE
$wsReco.Range("A20").Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$wsReco.Range("B20").Value = 0
$wsReco.Range("C20").Value = 4
$wsReco.Range("D20").Value = 419.05
$wsReco.Range("E20").Value = 105.57
$wsReco.Range("F20").Value = '🟡 Observer'
$wsReco.Range("G20").Value = '➖ Neutre'

# Row 21: BRVM - TELECOMMUNICATIONS
$wsReco.Range("A21").Value = 'BRVM - TELECOMMUNICATIONS'
$wsReco.Range("B21").Value = 0
$wsReco.Range("C21").Value = 4
$wsReco.Range("D21").Value = 368.2
$wsReco.Range("E21").Value = 92.06
$wsReco.Range("F21").Value = '🟡 Observer'
$wsReco.Range("G21").Value = '➖ Neutre'

# Row 22: BERNABE CI (BNBC)
$wsReco.Range("A22").Value = 'BERNABE CI (BNBC)'
$wsReco.Range("B22").Value = 2
$wsReco.Range("C22").Value = 0
$wsReco.Range("D22").Value = 12.66
$wsReco.Range("E22").Value = 3.33
$wsReco.Range("F22").Value = '🟡 Observer'
$wsReco.Range("G22").Value = '➖ Neutre'

# Row 23: VIVO ENERGY CI (SHEC)
$wsReco.Range("A23").Value = 'VIVO ENERGY CI (SHEC)'
$wsReco.Range("B23").Value = 2
$wsReco.Range("C23").Value = 0
$wsReco.Range("D23").Value = 10.82
$wsReco.Range("E23").Value = 5.73
$wsReco.Range("F23").Value = '🟡 Observer'
$wsReco.Range("G23").Value = '➖ Neutre'

# Row 24: SAFCA CI (SAFC)
$wsReco.Range("A24").Value = 'SAFCA CI (SAFC)'
$wsReco.Range("B24").Value = 3
$wsReco.Range("C24").Value = 1
$wsReco.Range("D24").Value = 10.3
$wsReco.Range("E24").Value = -7.2
$wsReco.Range("F24").Value = '🟢 Achat'
$wsReco.Range("G24").Value = '✅ Renforcer'

# Row 25: SETAO CI (STAC)
$wsReco.Range("A25").Value = 'SETAO CI (STAC)'
$wsReco.Range("B25").Value = 1
$wsReco.Range("C25").Value = 0
$wsReco.Range("D25").Value = 7.41
$wsReco.Range("E25").Value = 7.41
$wsReco.Range("F25").Value = '🟡 Observer'
$wsReco.Range("G25").Value = '➖ Neutre'

# Row 26: ECOBANK TRANS. INCORP. TG (ETIT)
$wsReco.Range("A26").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$wsReco.Range("B26").Value = 2
$wsReco.Range("C26").Value = 1
$wsReco.Range("D26").Value = 6.67
$wsReco.Range("E26").Value = 6.25
$wsReco.Range("F26").Value = '🟡 Observer'
$wsReco.Range("G26").Value = '👀 À surveiller'

# Row 27: ORAGROUP TOGO (ORGT)
$wsReco.Range("A27").Value = 'ORAGROUP TOGO (ORGT)'
$wsReco.Range("B27").Value = 2
$wsReco.Range("C27").Value = 1
$wsReco.Range("D27").Value = 5.63
$wsReco.Range("E27").Value = -0.6
$wsReco.Range("F27").Value = '🟡 Observer'
$wsReco.Range("G27").Value = '👀 À surveiller'

# Row 28: PALM CI (PALC)
$wsReco.Range("A28").Value = 'PALM CI (PALC)'
$wsReco.Range("B28").Value = 1
$wsReco.Range("C28").Value = 0
$wsReco.Range("D28").Value = 3.41
$wsReco.Range("E28").Value = 3.41
$wsReco.Range("F28").Value = '🟡 Observer'
$wsReco.Range("G28").Value = '➖ Neutre'

# Row 29: SITAB CI (STBC)
$wsReco.Range("A29").Value = 'SITAB CI (STBC)'
$wsReco.Range("B29").Value = 1
$wsReco.Range("C29").Value = 0
$wsReco.Range("D29").Value = 3.06
$wsReco.Range("E29").Value = 3.06
$wsReco.Range("F29").Value = '🟡 Observer'
$wsReco.Range("G29").Value = '➖ Neutre'

# Row 30: SERVAIR ABIDJAN CI (ABJC)
$wsReco.Range("A30").Value = 'SERVAIR ABIDJAN CI (ABJC)'
$wsReco.Range("B30").Value = 1
$wsReco.Range("C30").Value = 1
$wsReco.Range("D30").Value = 2.65
$wsReco.Range("E30").Value = -1.42
$wsReco.Range("F30").Value = '🟡 Observer'
$wsReco.Range("G30").Value = '👀 À surveiller'

# Row 31: SOCIETE GENERALE COTE D'IVOIRE (SGBC)
$wsReco.Range("A31").Value = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$wsReco.Range("B31").Value = 1
$wsReco.Range("C31").Value = 1
$wsReco.Range("D31").Value = 0.33
$wsReco.Range("E31").Value = 4.04
$wsReco.Range("F31").Value = '🟡 Observer'
$wsReco.Range("G31").Value = '👀 À surveiller'

# Row 32: UNIWAX CI (UNXC)
$wsReco.Range("A32").Value = 'UNIWAX CI (UNXC)'
$wsReco.Range("B32").Value = 1
$wsReco.Range("C32").Value = 1
$wsReco.Range("D32").Value = 0.28
$wsReco.Range("E32").Value = 5.45
$wsReco.Range("F32").Value = '🟡 Observer'
$wsReco.Range("G32").Value = '👀 À surveiller'

# Row 33: TOTAL
$wsReco.Range("A33").Value = 'TOTAL'
$wsReco.Range("B33").Value = 0
$wsReco.Range("C33").Value = 4
$wsReco.Range("D33").Value = 0
$wsReco.Range("E33").Value = 0
$wsReco.Range("F33").Value = '🟡 Observer'
$wsReco.Range("G33").Value = '➖ Neutre'

# Row 34: BICI CI (BICC)
$wsReco.Range("A34").Value = 'BICI CI (BICC)'
$wsReco.Range("B34").Value = 0
$wsReco.Range("C34").Value = 1
$wsReco.Range("D34").Value = -0.6
$wsReco.Range("E34").Value = -0.6
$wsReco.Range("F34").Value = '🟡 Observer'
$wsReco.Range("G34").Value = '➖ Neutre'

# Row 35: SAPH CI (SPHC)
$wsReco.Range("A35").Value = 'SAPH CI (SPHC)'
$wsReco.Range("B35").Value = 1
$wsReco.Range("C35").Value = 2
$wsReco.Range("D35").Value = -1.42
$wsReco.Range("E35").Value = 4.03
$wsReco.Range("F35").Value = '🟡 Observer'
$wsReco.Range("G35").Value = '👀 À surveiller'

# Row 36: AFRICA GLOBAL LOGISTICS CI (SDSC)
$wsReco.Range("A36").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$wsReco.Range("B36").Value = 0
$wsReco.Range("C36").Value = 1
$wsReco.Range("D36").Value = -1.75
$wsReco.Range("E36").Value = -1.75
$wsReco.Range("F36").Value = '🟡 Observer'
$wsReco.Range("G36").Value = '➖ Neutre'

# Row 37: FILTISAC CI (FTSC)
$wsReco.Range("A37").Value = 'FILTISAC CI (FTSC)'
$wsReco.Range("B37").Value = 0
$wsReco.Range("C37").Value = 1
$wsReco.Range("D37").Value = -1.84
$wsReco.Range("E37").Value = -1.84
$wsReco.Range("F37").Value = '🟡 Observer'
$wsReco.Range("G37").Value = '➖ Neutre'

# Row 38: TRACTAFRIC MOTORS CI (PRSC)
$wsReco.Range("A38").Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$wsReco.Range("B38").Value = 0
$wsReco.Range("C38").Value = 1
$wsReco.Range("D38").Value = -1.87
$wsReco.Range("E38").Value = -1.87
$wsReco.Range("F38").Value = '🟡 Observer'
$wsReco.Range("G38").Value = '➖ Neutre'

# Row 39: NSIA BANQUE COTE D'IVOIRE (NSBC)
$wsReco.Range("A39").Value = 'NSIA BANQUE COTE D''IVOIRE (NSBC)'
$wsReco.Range("B39").Value = 0
$wsReco.Range("C39").Value = 1
$wsReco.Range("D39").Value = -2
$wsReco.Range("E39").Value = -2
$wsReco.Range("F39").Value = '🟡 Observer'
$wsReco.Range("G39").Value = '➖ Neutre'

# Row 40: CIE CI (CIEC)
$wsReco.Range("A40").Value = 'CIE CI (CIEC)'
$wsReco.Range("B40").Value = 0
$wsReco.Range("C40").Value = 1
$wsReco.Range("D40").Value = -2.27
$wsReco.Range("E40").Value = -2.27
$wsReco.Range("F40").Value = '🟡 Observer'
$wsReco.Range("G40").Value = '➖ Neutre'

# Row 41: BANK OF AFRICA NG (BOAN)
$wsReco.Range("A41").Value = 'BANK OF AFRICA NG (BOAN)'
$wsReco.Range("B41").Value = 0
$wsReco.Range("C41").Value = 1
$wsReco.Range("D41").Value = -2.71
$wsReco.Range("E41").Value = -2.71
$wsReco.Range("F41").Value = '🟡 Observer'
$wsReco.Range("G41").Value = '➖ Neutre'

# Row 42: SOGB CI (SOGC)
$wsReco.Range("A42").Value = 'SOGB CI (SOGC)'
$wsReco.Range("B42").Value = 0
$wsReco.Range("C42").Value = 1
$wsReco.Range("D42").Value = -3.23
$wsReco.Range("E42").Value = -3.23
$wsReco.Range("F42").Value = '🟡 Observer'
$wsReco.Range("G42").Value = '➖ Neutre'

# Row 43: TOTALENERGIES MARKETING CI (TTLC)
$wsReco.Range("A43").Value = 'TOTALENERGIES MARKETING CI (TTLC)'
$wsReco.Range("B43").Value = 1
$wsReco.Range("C43").Value = 1
$wsReco.Range("D43").Value = -3.46
$wsReco.Range("E43").Value = 3.09
$wsReco.Range("F43").Value = '🟡 Observer'
$wsReco.Range("G43").Value = '👀 À surveiller'

# Row 44: SOLIBRA CI (SLBC)
$wsReco.Range("A44").Value = 'SOLIBRA CI (SLBC)'
$wsReco.Range("B44").Value = 0
$wsReco.Range("C44").Value = 2
$wsReco.Range("D44").Value = -5.2
$wsReco.Range("E44").Value = -1.97
$wsReco.Range("F44").Value = '🟡 Observer'
$wsReco.Range("G44").Value = '➖ Neutre'

# --- Sheet "Top_YTD": refresh YTD progression figures ---
# Row 2: BRVM - SERVICES PUBLICS
$wsYtd.Range("A2").Value = 'BRVM - SERVICES PUBLICS'
$wsYtd.Range("B2").Value = 6747097.64

# Row 3: BRVM - AUTRES SECTEURS
$wsYtd.Range("A3").Value = 'BRVM - AUTRES SECTEURS'
$wsYtd.Range("B3").Value = 296648.23

# Row 4: NEI-CEDA CI
$wsYtd.Range("A4").Value = 'NEI-CEDA CI'
$wsYtd.Range("B4").Value = 234891.66

# Row 5: UNIWAX CI
$wsYtd.Range("A5").Value = 'UNIWAX CI'
$wsYtd.Range("B5").Value = 204269.75

# Row 6: AIR LIQUIDE CI
$wsYtd.Range("A6").Value = 'AIR LIQUIDE CI'
$wsYtd.Range("B6").Value = 172956

# Row 7: SETAO CI
$wsYtd.Range("A7").Value = 'SETAO CI'
$wsYtd.Range("B7").Value = 168881.12

# Row 8: BRVM - DISTRIBUTION
$wsYtd.Range("A8").Value = 'BRVM - DISTRIBUTION'
$wsYtd.Range("B8").Value = 46488.86

# Row 9: BRVM - TRANSPORT
$wsYtd.Range("A9").Value = 'BRVM - TRANSPORT'
$wsYtd.Range("B9").Value = 40467.57

# Row 10: CFAO MOTORS CI
$wsYtd.Range("A10").Value = 'CFAO MOTORS CI'
$wsYtd.Range("B10").Value = 40141.25

# Row 11: BRVM - AGRICULTURE
$wsYtd.Range("A11").Value = 'BRVM - AGRICULTURE'
$wsYtd.Range("B11").Value = 31533.07

